$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-02-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-17 Monday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "24+24=48"
$t.Cell(1,2).Range.Text = "3+40=43"
$t.Cell(1,3).Range.Text = "56+11=67"
$t.Cell(1,4).Range.Text = "38+53=91"
$t.Cell(1,5).Range.Text = "27+33=60"
$t.Cell(2,1).Range.Text = "96-82=14"
$t.Cell(2,2).Range.Text = "32+30=62"
$t.Cell(2,3).Range.Text = "44-16=28"
$t.Cell(2,4).Range.Text = "45+25=70"
$t.Cell(2,5).Range.Text = "60-4=56"
$t.Cell(3,1).Range.Text = "19-3=16"
$t.Cell(3,2).Range.Text = "74-58=16"
$t.Cell(3,3).Range.Text = "6+70=76"
$t.Cell(3,4).Range.Text = "34-20=14"
$t.Cell(3,5).Range.Text = "55-10=45"
$t.Cell(4,1).Range.Text = "82-68=14"
$t.Cell(4,2).Range.Text = "80-57=23"
$t.Cell(4,3).Range.Text = "17-13=4"
$t.Cell(4,4).Range.Text = "18+6=24"
$t.Cell(4,5).Range.Text = "97-17=80"
$t.Cell(5,1).Range.Text = "28-18=10"
$t.Cell(5,2).Range.Text = "1+50=51"
$t.Cell(5,3).Range.Text = "88-80=8"
$t.Cell(5,4).Range.Text = "13+31=44"
$t.Cell(5,5).Range.Text = "80-29=51"
$t.Cell(6,1).Range.Text = "77-19=58"
$t.Cell(6,2).Range.Text = "89-87=2"
$t.Cell(6,3).Range.Text = "1+61=62"
$t.Cell(6,4).Range.Text = "27+43=70"
$t.Cell(6,5).Range.Text = "70-42=28"
$t.Cell(7,1).Range.Text = "30-18=12"
$t.Cell(7,2).Range.Text = "58+37=95"
$t.Cell(7,3).Range.Text = "8+46=54"
$t.Cell(7,4).Range.Text = "12+32=44"
$t.Cell(7,5).Range.Text = "99-86=13"
$t.Cell(8,1).Range.Text = "83-35=48"
$t.Cell(8,2).Range.Text = "6+53=59"
$t.Cell(8,3).Range.Text = "75-66=9"
$t.Cell(8,4).Range.Text = "27+10=37"
$t.Cell(8,5).Range.Text = "52+38=90"
$t.Cell(9,1).Range.Text = "21+4=25"
$t.Cell(9,2).Range.Text = "85-69=16"
$t.Cell(9,3).Range.Text = "61+14=75"
$t.Cell(9,4).Range.Text = "33-2=31"
$t.Cell(9,5).Range.Text = "31+53=84"
$t.Cell(10,1).Range.Text = "43+6=49"
$t.Cell(10,2).Range.Text = "55-0=55"
$t.Cell(10,3).Range.Text = "8+64=72"
$t.Cell(10,4).Range.Text = "11+36=47"
$t.Cell(10,5).Range.Text = "42-36=6"
$t.Cell(11,1).Range.Text = "94-2=92"
$t.Cell(11,2).Range.Text = "89-54=35"
$t.Cell(11,3).Range.Text = "68-20=48"
$t.Cell(11,4).Range.Text = "62-61=1"
$t.Cell(11,5).Range.Text = "58-23=35"
$t.Cell(12,1).Range.Text = "3+23=26"
$t.Cell(12,2).Range.Text = "67-58=9"
$t.Cell(12,3).Range.Text = "87+6=93"
$t.Cell(12,4).Range.Text = "27+67=94"
$t.Cell(12,5).Range.Text = "83-71=12"
$t.Cell(13,1).Range.Text = "68-40=28"
$t.Cell(13,2).Range.Text = "69-58=11"
$t.Cell(13,3).Range.Text = "99-60=39"
$t.Cell(13,4).Range.Text = "16+82=98"
$t.Cell(13,5).Range.Text = "94-62=32"
$t.Cell(14,1).Range.Text = "86-14=72"
$t.Cell(14,2).Range.Text = "67-58=9"
$t.Cell(14,3).Range.Text = "55+11=66"
$t.Cell(14,4).Range.Text = "39-26=13"
$t.Cell(14,5).Range.Text = "54+21=75"
$t.Cell(15,1).Range.Text = "10+86=96"
$t.Cell(15,2).Range.Text = "95-74=21"
$t.Cell(15,3).Range.Text = "61+14=75"
$t.Cell(15,4).Range.Text = "68-53=15"
$t.Cell(15,5).Range.Text = "73+22=95"
$t.Cell(16,1).Range.Text = "97-49=48"
$t.Cell(16,2).Range.Text = "68-56=12"
$t.Cell(16,3).Range.Text = "9-8=1"
$t.Cell(16,4).Range.Text = "61-55=6"
$t.Cell(16,5).Range.Text = "10+3=13"
$t.Cell(17,1).Range.Text = "11+75=86"
$t.Cell(17,2).Range.Text = "37-8=29"
$t.Cell(17,3).Range.Text = "7+38=45"
$t.Cell(17,4).Range.Text = "81-9=72"
$t.Cell(17,5).Range.Text = "41+10=51"
$t.Cell(18,1).Range.Text = "61+20=81"
$t.Cell(18,2).Range.Text = "71-17=54"
$t.Cell(18,3).Range.Text = "19+69=88"
$t.Cell(18,4).Range.Text = "71-24=47"
$t.Cell(18,5).Range.Text = "77-54=23"
$t.Cell(19,1).Range.Text = "4+1=5"
$t.Cell(19,2).Range.Text = "32+1=33"
$t.Cell(19,3).Range.Text = "18+55=73"
$t.Cell(19,4).Range.Text = "45+5=50"
$t.Cell(19,5).Range.Text = "21+8=29"
$t.Cell(20,1).Range.Text = "70+22=92"
$t.Cell(20,2).Range.Text = "89-20=69"
$t.Cell(20,3).Range.Text = "74-5=69"
$t.Cell(20,4).Range.Text = "63-61=2"
$t.Cell(20,5).Range.Text = "13+16=29"
